$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{ A = 44; B = '41'; C = 'Access Control'; D = 'Constant across alternatives' },
  @{ A = 45; B = '41'; C = 'Closures'; D = '4 temporal closure alternatives.  Base case, add May closure, add October closure, closure rule - if sample below threshold catch close' },
  @{ A = 46; B = '42'; C = 'Catch Limit'; D = 'Two HCRs, based on a static or a dynamic B(0) used as a BRP.' },
  @{ A = 47; B = '43'; C = 'Catch Limit'; D = 'Two HCRs, one with a constant F, and one with a variable F based on BRPs' },
  @{ A = 48; B = '43'; C = 'Catch Limit'; D = 'Two additional HCRs, which correspond to the first two, but are adjusted based on environmental conditions.' },
  @{ A = 49; B = '44'; C = 'Catch Limit'; D = '4 management procedures covering different approaches to developing a rebuilding plan based on the acceptable recovery probability or inclusion of climate effects in the plan.' },
  @{ A = 50; B = '45'; C = 'Catch Limit'; D = 'This is a TAC based fishery.  The TACs weren''t the focus, but would be altered by reference point changes.' },
  @{ A = 51; B = '45'; C = 'Other'; D = 'Two reference point options, one that shifts to account for climate change, and one that does not.' },
  @{ A = 52; B = '46'; C = 'Catch Limit'; D = 'TAC implemented or not' },
  @{ A = 53; B = '46'; C = 'Size Limit'; D = 'Size limit implemented or not' },
  @{ A = 54; B = '46'; C = 'Closures'; D = 'marine protected areas implemented or not' },
  @{ A = 55; B = '46'; C = 'Other'; D = 'Pollution reduction implemented or not' },
  @{ A = 56; B = '47'; C = 'Other'; D = 'Note, not applicable to fishery management.  There are three approaches to beach replenishment, fixed amount and interval, fixed amount, and fixed interval, as well as no action.' },
  @{ A = 57; B = '48'; C = 'Closures'; D = 'Not cleared provided.' },
  @{ A = 58; B = '31'; C = 'Closures'; D = 'Spatial closures: No closures, 3 closure durations, 2 closure location rules' },
  @{ A = 59; B = '31'; C = 'Other'; D = 'Size based closure rules: 4 options' },
  @{ A = 60; B = '49'; C = 'Closures'; D = 'Spatial closures: No closures, 3 closure durations, 2 closure location rules' },
  @{ A = 61; B = '49'; C = 'Other'; D = 'Size based closure rules: 4 options' },
  @{ A = 62; B = '50'; C = 'Catch Limit'; D = '6 alternative methods for adjusting the catch limit.' },
  @{ A = 63; B = '32'; C = 'Access Control, Catch Limit'; D = 'Quota system, accompanied by location and gear restrictions, combined into 4 strategy' },
  @{ A = 64; B = '33'; C = 'Size Limit'; D = '# of size limits utilized and areas to which they are applied' },
  @{ A = 65; B = '34'; C = 'Catch Limit, Effort Limit'; D = '22 MPs available in the DLMtoolkit package in R, 11 output control MPs and 11 input control MPs' },
  @{ A = 66; B = '35'; C = 'Catch Limit'; D = 'TAC=ABC set based on ICES F based proceedures (with and without uncertainty) and a constant F' },
  @{ A = 67; B = '36'; C = 'Catch Limit'; D = '4 HCRs.  2 based on BRPs and 2 proportional harvest rules, the difference in each category is whether there is an annual TAC change limit.' },
  @{ A = 68; B = '36'; C = 'Other'; D = '3 stock assessment methods.  XSA, Schaefer, Difference' },
  @{ A = 69; B = '37'; C = 'Effort Limit'; D = '2 levels, Status quo and reduction from 9000 to 5000 boat days' },
  @{ A = 70; B = '37'; C = 'Closures'; D = '4 options, status quo, reef buffer area closure, Masig area closure, moon cycle calendar closure' },
  @{ A = 71; B = '38'; C = 'Catch Limit'; D = 'No alternatives, a single Harvest quota accompanied by a season duration' },
  @{ A = 72; B = '39'; C = 'Catch Limit'; D = '4 TACs; 750 t, 1000 t, 1250 t, and 1500 t; competitave TAC vs ITQ' },
  @{ A = 73; B = '39'; C = 'Effort Limit'; D = 'constant effort limit' },
  @{ A = 74; B = '39'; C = 'Closures'; D = 'four spatially explicit no-take extents: an extent consistent with that from the mid-1980s to mid-2004 (approximately 16% of coral trout habitat in the park); an extent implemented during rezoning in 2004 (32%); a hypothetical extent of 50% (Little et al. 2009a); and a hypothetical extent of 0%.' },
  @{ A = 75; B = '40'; C = 'Catch Limit'; D = '6 HCRs combining timeline and precaution:  a reactive decision interval with no additional ACL reduction, and five HCRs consisting of a fixed decision interval with precautionary ACL reductions of 0 (i.e., no reduction), 10, 20, 30, and 40%.' }
)

$startRow = 44
$srcRow = 43
$i = 0
foreach ($item in $rows) {
  $r = $startRow + $i

  $aCell = $ws.Cells.Item($r, 1)
  $aCell.Value = $item.A
  $ws.Range("A$srcRow").Copy()
  $aCell.PasteSpecial(-4122)

  $bCell = $ws.Cells.Item($r, 2)
  $bCell.NumberFormat = "@"
  $bCell.Value = $item.B
  $ws.Range("B$srcRow").Copy()
  $bCell.PasteSpecial(-4122)

  $cCell = $ws.Cells.Item($r, 3)
  $cCell.Value = $item.C
  $ws.Range("C$srcRow").Copy()
  $cCell.PasteSpecial(-4122)

  $dCell = $ws.Cells.Item($r, 4)
  $dCell.Value = $item.D
  $ws.Range("D$srcRow").Copy()
  $dCell.PasteSpecial(-4122)

  $i++
}

$excel.CutCopyMode = 0
